$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25 (bus index 0-23), columns B,C,D,E,F,I,J,K,L,M,N.
# Column G (slack bus, always 1.0) and column H (unused) are left unchanged.
$newData = @{
    2 = @(1.02, 1.018126368627203, 1.033750616651499, 1.019437735780974, 1.016500239982106, 1.033004278808795, 1.023336572893786, 1.036552066835184, 1.022280984482213, 1.01935222567328, 1.011878424093096)
    3 = @(1.02, 1.019002730153416, 1.034281418323099, 1.020178580846987, 1.018016794697574, 1.033182753010039, 1.023849504065007, 1.036892656286116, 1.022828033653078, 1.020672206734843, 1.0120484931492)
    4 = @(1.02, 1.019569721580501, 1.034624434110795, 1.020658300281357, 1.018997990295555, 1.033296583308165, 1.024180721593766, 1.037111882458369, 1.023181691898176, 1.021525725634811, 1.012158286931183)
    5 = @(1.02, 1.019808066971531, 1.034768529018679, 1.020860055723486, 1.019410461683252, 1.033344040912284, 1.024319801631232, 1.037203766845779, 1.023330293002031, 1.021884405986954, 1.012204383725038)
    6 = @(1.02, 1.019848085158283, 1.034792716753799, 1.020893936100722, 1.019479716217703, 1.033351985976054, 1.024343144168324, 1.03721917828798, 1.023355239264082, 1.021944621960519, 1.012212120027622)
    7 = @(1.02, 1.019572906433641, 1.03462635994404, 1.020660995828602, 1.019003501843909, 1.033297218998303, 1.024182580632023, 1.037113111316928, 1.023183677815798, 1.021530518880472, 1.012158903116724)
    8 = @(1.02, 1.018422554056181, 1.033930095454196, 1.019688036338805, 1.017012793746039, 1.033064937219851, 1.023510061709567, 1.036667409758931, 1.022465928143746, 1.019798445472985, 1.011935951890621)
    9 = @(1.02, 1.016394939374603, 1.032699823777113, 1.017976218183439, 1.01350380723886, 1.032642981026125, 1.022319774740211, 1.035873204870938, 1.021198735064436, 1.01674156023991, 1.011541156265519)
    10 = @(1.02, 1.015042843025555, 1.031877473104138, 1.016836837199394, 1.011163444277682, 1.03235320717016, 1.021522750041365, 1.035337865934529, 1.020352330381798, 1.014700192941777, 1.011276669561317)
    11 = @(1.02, 1.014457288683784, 1.031520890858011, 1.016343915235058, 1.010149728582304, 1.03222572955288, 1.02117680107701, 1.03510467902914, 1.019985450147218, 1.013815385425596, 1.011161839062598)
    12 = @(1.02, 1.014239774612262, 1.031388366893286, 1.016160888404526, 1.009773135108651, 1.032178078144166, 1.021048175279749, 1.035017856448413, 1.019849117442928, 1.013486591092259, 1.011119139984594)
    13 = @(1.02, 1.014286432711688, 1.031416797032638, 1.016200145269645, 1.009853918238817, 1.032188313125832, 1.021075771603331, 1.035036489519148, 1.019878363848415, 1.013557124890822, 1.011128301158775)
    14 = @(1.02, 1.014439309170051, 1.031509937875167, 1.016328784823504, 1.010118600390151, 1.032221796800878, 1.021166171377125, 1.035097506456038, 1.019974182005708, 1.013788210015382, 1.011158310479884)
    15 = @(1.02, 1.014533499665089, 1.031567315335912, 1.016408052706061, 1.010281672426548, 1.0322423873871, 1.021221853105846, 1.035135073660084, 1.020033211182634, 1.013930570829765, 1.011176794130125)
    16 = @(1.02, 1.015081702496184, 1.031901127910938, 1.016869560076421, 1.011230713949184, 1.032361625266491, 1.02154569200986, 1.035353312757208, 1.020376671000465, 1.014758895700761, 1.011284284040944)
    17 = @(1.02, 1.01542555221303, 1.032110387209546, 1.017159169021778, 1.011825932045927, 1.032435883948799, 1.021748604884998, 1.035489839182159, 1.020592012335036, 1.015278242409156, 1.011351627735422)
    18 = @(1.02, 1.015626105470093, 1.0322323963724, 1.017328135194392, 1.012173081907604, 1.032479004429416, 1.021866880252121, 1.035569339364465, 1.020717580473944, 1.015581083787574, 1.011390878643835)
    19 = @(1.02, 1.015694487510222, 1.032273990114074, 1.017385755427854, 1.012291445993514, 1.032493674588468, 1.021907195477239, 1.035596424219764, 1.020760389711459, 1.015684330642225, 1.011404257178835)
    20 = @(1.02, 1.015388661267616, 1.032087940646428, 1.017128092356013, 1.011762074002515, 1.03242793669082, 1.021726842557091, 1.035475204976872, 1.020568912054157, 1.015222530225474, 1.011344405448239)
    21 = @(1.02, 1.014394291225071, 1.03148251223369, 1.01629090184521, 1.010040659639959, 1.032211944985872, 1.021139554345713, 1.035079544183683, 1.01994596751944, 1.013720165039043, 1.011149474749899)
    22 = @(1.02, 1.013769015444738, 1.0311014311953, 1.015764911024117, 1.008958020894241, 1.032074403795034, 1.020769580044018, 1.034829581806211, 1.019553967613325, 1.012774770865377, 1.01102664844196)
    23 = @(1.02, 1.014100493140297, 1.031303489164581, 1.016043712105671, 1.009531980185711, 1.032147481627509, 1.020965778928385, 1.034962204510688, 1.019761805373327, 1.013276019450336, 1.011091786153306)
    24 = @(1.02, 1.015405330723922, 1.032098083431323, 1.017142134433265, 1.011790928797342, 1.03243152831205, 1.021736676264274, 1.035481817956493, 1.020579350190447, 1.015247704423911, 1.01134766898015)
    25 = @(1.02, 1.016919189657705, 1.033018267631323, 1.018418444512464, 1.014411125661084, 1.032753562038509, 1.022628110165565, 1.03607956484805, 1.021526620088533, 1.017532426973037, 1.011643448133758)
}

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

foreach ($r in $newData.Keys) {
    $rowValues = $newData[$r]
    for ($i = 0; $i -lt $cols.Count; $i++) {
        $cellRef = "$($cols[$i])$r"
        $ws.Range($cellRef).Value = $rowValues[$i]
    }
}